$d = $word.ActiveDocument

# Remove the stray apostrophe after "asterisks" (asterisks' -> asterisks)
$d.Content.Find.Execute("asterisks’ will", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "asterisks will", 2)

# Re-locate the middle portion of the sentence and toggle (then untoggle) a
# character-level formatting property on it. Word splits the containing run
# at the boundaries of the re-formatted range; toggling Bold back to its
# original value leaves the visible/semantic formatting untouched (no
# leftover <w:b/> in the run properties) while keeping the run split, which
# yields the desired three runs:
#   "The feature of " | "hiding passwords with asterisks" | " will be done with ncursors."
$r = $d.Content
$r.Find.Execute("hiding passwords with asterisks", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
$r.Bold = $true
$r.Bold = $false
